$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 06:23"

# Row 5 - India
$ws.Range("B5").Value = 4280422
$ws.Range("C5").Value = 2838
$ws.Range("D5").Value = 3323950
$ws.Range("E5").Value = 883656

# Row 20 - Pakistan
$ws.Range("B20").Value = 299233
$ws.Range("C20").Value = 330
$ws.Range("D20").Value = 286157
$ws.Range("E20").Value = 6726
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 6350

# Row 50 - Honduras
$ws.Range("B50").Value = 64814
$ws.Range("C50").Value = 50
$ws.Range("D50").Value = 13828
$ws.Range("E50").Value = 48963
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 2023

# Row 125 - Tailandia
$ws.Range("B125").Value = 3446
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 3284
$ws.Range("E125").Value = 104

# Row 187 - Butan
$ws.Range("B187").Value = 233
$ws.Range("C187").Value = 3
$ws.Range("E187").Value = 82
